$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z3").Value = "2025-11-13T06:53:00.395831"
$ws.Range("Z4:Z13").Value = "2025-11-13T06:53:00.396833"
$ws.Range("Z14:Z24").Value = "2025-11-13T06:53:00.397833"
$ws.Range("Z25:Z32").Value = "2025-11-13T06:53:00.398832"
$ws.Range("Z33:Z34").Value = "2025-11-13T06:53:00.401240"
$ws.Range("Z35:Z45").Value = "2025-11-13T06:53:00.401758"
$ws.Range("Z46:Z50").Value = "2025-11-13T06:53:00.625315"
$ws.Range("Z51:Z54").Value = "2025-11-13T06:53:00.626315"
$ws.Range("Z55:Z59").Value = "2025-11-13T06:53:00.627315"
$ws.Range("Z60:Z64").Value = "2025-11-13T06:53:00.628315"
$ws.Range("Z65:Z68").Value = "2025-11-13T06:53:00.629315"
$ws.Range("Z69:Z70").Value = "2025-11-13T06:53:00.630315"
$ws.Range("Z71:Z74").Value = "2025-11-13T06:53:00.631316"
$ws.Range("Z75:Z77").Value = "2025-11-13T06:53:01.094335"
$ws.Range("Z78:Z83").Value = "2025-11-13T06:53:01.095336"
$ws.Range("Z84:Z85").Value = "2025-11-13T06:53:01.096336"
$ws.Range("Z86:Z87").Value = "2025-11-13T06:53:01.096862"
$ws.Range("Z88:Z91").Value = "2025-11-13T06:53:01.097399"
$ws.Range("Z92").Value = "2025-11-13T06:53:01.098390"
$ws.Range("Z93").Value = "2025-11-13T06:53:01.099391"
$ws.Range("Z94").Value = "2025-11-13T06:53:01.101954"
$ws.Range("Z95").Value = "2025-11-13T06:53:01.102953"
$ws.Range("Z96:Z102").Value = "2025-11-13T06:53:01.110195"

Write-Output "done"
